# ws.async-showcase.xlsx maintenance edit
#
# [json] add new function `storeKeys(json,jsonpath,var)` to the `json`
#        command list on the hidden `#system` sheet (inserted in
#        alphabetical order, right before `storeValue`).
#
# Also retires the now-unused `text` category (its single command,
# `spellCheck(var,profile,text)`, has no callers) which frees up column Y
# and removes the `text` entry from the `target` category list; the
# categories that used to live to the right of it (web, webalert,
# webcookie, ws, ws.async, xml) shift one row/column to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1) json: insert "storeKeys(json,jsonpath,var)" before "storeValue" ---
# Column M currently holds: ... M15=storeCount, M16=storeValue, M17=storeValues
# Shift M16:M17 down into M17:M18 (bottom-up so we don't clobber data we
# still need to read), then drop the new entry into the now-free M16.
$m17 = $ws.Range("M17").Value()
$m16 = $ws.Range("M16").Value()
$ws.Range("M18").Value = $m17
$ws.Range("M17").Value = $m16
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# --- 2) target: drop the "text" row (A25), shifting the rest up ---
$ws.Range("A25").Delete()

# --- 3) text: remove the now-orphaned column entirely; everything right
#        of it (web, webalert, webcookie, ws, ws.async, xml) shifts left
$ws.Columns("Y").Delete()

# --- 4) fix up the named ranges that moved as a result of the above ---
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
